$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "comments" column header
$ws.Range("F1").Value = "comments"

# Add comment values for each row
$ws.Range("F2").Value = "AAA"
$ws.Range("F3").Value = "BBB"
$ws.Range("F4").Value = "CCC"
$ws.Range("F5").Value = "DDD"

# Update price column values (B3, B4, B5)
$ws.Range("B3").Value = 6.7565
$ws.Range("B4").Value = 4.364
$ws.Range("B5").Value = 5.42

# Update rating column values (E4, E5)
$ws.Range("E4").Value = 3
$ws.Range("E5").Value = 65.1234

# Remove custom currency number format on price column (B2:B5) -> General format
$ws.Range("B2:B5").NumberFormat = "General"

# Update selection to F6
$ws.Range("F6").Select()
